$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    heading ("Play Dollar Bomb Free: Slot Game Review"). It consists
#    of a leading empty run, a bold "Meta description" run and a
#    normal run with the rest of the description text, matching the
#    structure used elsewhere in the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Discover more about Dollar Bomb by CQ9 Gaming in our review. Play Dollar Bomb for free online and experience great winning potential with an eastern theme.</w:t></w:r></w:p>'
$metaPara.Range.InsertXML($metaXml) | Out-Null

# ------------------------------------------------------------------
# 2) Near the end of the document there is a duplicated
#    "Play Dollar Bomb Free: Slot Game Review" paragraph (bold) right
#    before the closing italic blurb paragraph. That duplicate is no
#    longer needed now that the title/meta-description live at the
#    top of the document, so remove it entirely. (The real title
#    paragraph at the very top uses the "Heading 1" style, while the
#    duplicate near the bottom is plain/body style with direct bold
#    character formatting, so checking the style keeps us from ever
#    touching the real title.)
# ------------------------------------------------------------------
$dupText = "Play Dollar Bomb Free: Slot Game Review"
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $dupText -and $p.Style.NameLocal -notlike "Heading*") {
        $p.Range.Delete() | Out-Null
        break
    }
}

# ------------------------------------------------------------------
# 3) Replace the closing italic blurb's text with new feature-image
#    generation instructions, preserving the existing italic run
#    formatting already applied to that paragraph.
# ------------------------------------------------------------------
$oldBlurb = "Discover more about Dollar Bomb by CQ9 Gaming in our review. Play Dollar Bomb for free online and experience great winning potential with an eastern theme."
$newBlurb = "Create a feature image for Dollar Bomb with the following specifications: Design a cartoon-style image featuring a happy Maya warrior with glasses. The warrior should be holding a bag of golden coins and standing in front of a grand castle. The castle should be in the background, with Dollar Bomb's name prominently displayed above it in bold, golden letters. Use bright, eye-catching colors to make the image stand out and add a touch of whimsy to the design. The image should convey excitement, joy, and the promise of big wins to players."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $oldBlurb) {
        $r = $p.Range
        $r.MoveEnd(1, -1) | Out-Null
        $r.Text = $newBlurb
        break
    }
}
